$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '57.487.69'
$ws.Cells.Item(2, 5).Value = '  -6.34%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.904.35'
$ws.Cells.Item(3, 5).Value = '  -4.06%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.07%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '551.43'
$ws.Cells.Item(5, 5).Value = '  -2.83%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '123.35'
$ws.Cells.Item(6, 5).Value = '  -4.97%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.10%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '2.900.22'
$ws.Cells.Item(8, 5).Value = '  -4.18%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.53%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -8.68%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '4.71'
$ws.Cells.Item(11, 5).Value = '  -11.08%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.95%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -5.92%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '32.61'
$ws.Cells.Item(14, 5).Value = '  -2.07%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +0.51%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.378.40'
$ws.Cells.Item(16, 5).Value = '  -4.15%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.894.98'
$ws.Cells.Item(17, 5).Value = '  -4.06%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '6.54'
$ws.Cells.Item(18, 5).Value = '  +4.65%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '57.462.79'
$ws.Cells.Item(19, 5).Value = '  -6.49%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '404.15'
$ws.Cells.Item(20, 5).Value = '  -8.36%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '12.91'
$ws.Cells.Item(21, 5).Value = '  -2.59%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.672'
$ws.Cells.Item(22, 5).Value = '  +0.66%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -4.84%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '12.81'
$ws.Cells.Item(24, 5).Value = '  -1.40%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '77.17'
$ws.Cells.Item(25, 5).Value = '  -2.88%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.04%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.998'
$ws.Cells.Item(27, 5).Value = '  -0.15%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -2.14%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +1.82%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.19'
$ws.Cells.Item(30, 5).Value = '  -1.31%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.03'
$ws.Cells.Item(31, 5).Value = '  -2.83%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '24.75'
$ws.Cells.Item(32, 5).Value = '  -3.51%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.0989'
$ws.Cells.Item(33, 5).Value = '  +4.62%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '5.43'
$ws.Cells.Item(34, 5).Value = '  -3.65%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.904'
$ws.Cells.Item(35, 5).Value = '  -6.15%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.99'
$ws.Cells.Item(36, 5).Value = '  -13.41%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '47.94'
$ws.Cells.Item(37, 5).Value = '  -4.91%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +6.57%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0₃0620'
$ws.Cells.Item(39, 5).Value = '  -8.49%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.106'
$ws.Cells.Item(40, 5).Value = '  -2.86%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'VeChain'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0341'
$ws.Cells.Item(41, 5).Value = '  -6.04%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.621.55'
$ws.Cells.Item(42, 5).Value = '  -2.49%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Bittensor'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '359.78'
$ws.Cells.Item(43, 5).Value = '  -5.90%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'dogwifhat'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.40'
$ws.Cells.Item(44, 5).Value = '  -3.06%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '118.91'
$ws.Cells.Item(46, 5).Value = '  -1.04%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -3.66%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -0.33%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.95'
$ws.Cells.Item(49, 5).Value = '  -2.44%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '22.85'
$ws.Cells.Item(50, 5).Value = '  -2.81%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -4.47%  '
